$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 56691736  # was 6084806
$ws.Range("I40").Value = 29425.5  # was 25040.4
$ws.Range("J40").Value = 94466610  # was 11134610
$ws.Range("K40").Value = 29425.5  # was 25040.4
$ws.Range("L40").Value = 94466610  # was 11134610
$ws.Range("M40").Value = -29250.5  # was -24865.4
$ws.Range("N40").Value = -94466960  # was -11134960

# Row 51
$ws.Range("H51").Value = 5970.737  # was 7066
$ws.Range("I51").Value = 4477.1113  # was 1494
$ws.Range("J51").Value = 7315  # was 7762.5
$ws.Range("K51").Value = 4477.1113  # was 1494
$ws.Range("L51").Value = 7315  # was 7762.5
$ws.Range("M51").Value = -3993.1113  # was -1010
$ws.Range("N51").Value = -8283  # was -8730.5

# Row 62
$ws.Range("H62").Value = 4363.6665  # was 4265.8
$ws.Range("I62").Value = 3526.6667  # was 3491.25
$ws.Range("K62").Value = 3526.6667  # was 3491.25
$ws.Range("M62").Value = -2902.6667  # was -2867.25

# Row 65
$ws.Range("H65").Value = 4363.6665  # was 4265.8
$ws.Range("I65").Value = 3526.6667  # was 3491.25
$ws.Range("K65").Value = 17633.3335  # was 17456.25
$ws.Range("M65").Value = -14513.3335  # was -14336.25

# Row 86
$ws.Range("H86").Value = 258079520  # was 172053340
$ws.Range("I86").Value = 344096030  # was 258072270
$ws.Range("J86").Value = 30004  # was 15504
$ws.Range("K86").Value = 344096030  # was 258072270
$ws.Range("L86").Value = 30004  # was 15504
$ws.Range("M86").Value = -344094907  # was -258071147
$ws.Range("N86").Value = -32250  # was -17750

# Row 89
$ws.Range("H89").Value = 258079520  # was 172053340
$ws.Range("I89").Value = 344096030  # was 258072270
$ws.Range("J89").Value = 30004  # was 15504
$ws.Range("K89").Value = 1720480150  # was 1290361350
$ws.Range("L89").Value = 150020  # was 77520
$ws.Range("M89").Value = -1720474534  # was -1290355734
$ws.Range("N89").Value = -161252  # was -88752

# Row 132
$ws.Range("H132").Value = 117976.75  # was 117980.42
$ws.Range("I132").Value = 282896.03  # was 282905.56
$ws.Range("K132").Value = 848688.0900000001  # was 848716.6799999999
$ws.Range("M132").Value = -846158.0900000001  # was -846186.6799999999

# Row 137
$ws.Range("H137").Value = 4782.4287  # was 4754.421
$ws.Range("J137").Value = 5880.077  # was 6031.273
$ws.Range("L137").Value = 17640.231  # was 18093.819
$ws.Range("N137").Value = -22740.231  # was -23193.819

# Row 140
$ws.Range("H140").Value = 60516.668  # was 60285
$ws.Range("J140").Value = 59331.25  # was 59205.555
$ws.Range("L140").Value = 59331.25  # was 59205.555
$ws.Range("N140").Value = -69691.25  # was -69565.55499999999


# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 728911.9  # was 699812.5
$ws.Range("I2").Value = 1027526.44  # was 1027533.9
$ws.Range("J2").Value = 3705.1428  # was 3404.5
$ws.Range("K2").Value = 1027526.44  # was 1027533.9
$ws.Range("L2").Value = 3705.1428  # was 3404.5
$ws.Range("M2").Value = -1027413.44  # was -1027420.9
$ws.Range("N2").Value = -3931.1428  # was -3630.5

# Row 32
$ws.Range("H32").Value = 2096.8206  # was 2070.9114
$ws.Range("I32").Value = 1675.6857  # was 1652.7887
$ws.Range("K32").Value = 1675.6857  # was 1652.7887
$ws.Range("M32").Value = -1388.6857  # was -1365.7887

# Row 74
$ws.Range("H74").Value = 4758.927  # was 4755.619
$ws.Range("I74").Value = 866.5  # was 1470.1666
$ws.Range("J74").Value = 5426.2  # was 5303.1943
$ws.Range("K74").Value = 866.5  # was 1470.1666
$ws.Range("L74").Value = 5426.2  # was 5303.1943
$ws.Range("M74").Value = 7.5  # was -596.1666
$ws.Range("N74").Value = -7174.2  # was -7051.1943

# Row 77
$ws.Range("H77").Value = 4758.927  # was 4755.619
$ws.Range("I77").Value = 866.5  # was 1470.1666
$ws.Range("J77").Value = 5426.2  # was 5303.1943
$ws.Range("K77").Value = 4332.5  # was 7350.833000000001
$ws.Range("L77").Value = 27131  # was 26515.9715
$ws.Range("M77").Value = 35.5  # was -2982.833000000001
$ws.Range("N77").Value = -35867  # was -35251.9715

# Row 116
$ws.Range("H116").Value = 728911.9  # was 699812.5
$ws.Range("I116").Value = 1027526.44  # was 1027533.9
$ws.Range("J116").Value = 3705.1428  # was 3404.5
$ws.Range("K116").Value = 1027526.44  # was 1027533.9
$ws.Range("L116").Value = 3705.1428  # was 3404.5
$ws.Range("M116").Value = -1025232.44  # was -1025239.9
$ws.Range("N116").Value = -8293.1428  # was -7992.5


# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 728911.9  # was 699812.5
$ws.Range("I3").Value = 1027526.44  # was 1027533.9
$ws.Range("J3").Value = 3705.1428  # was 3404.5
$ws.Range("K3").Value = 1027526.44  # was 1027533.9
$ws.Range("L3").Value = 3705.1428  # was 3404.5
$ws.Range("M3").Value = -1027412.44  # was -1027419.9
$ws.Range("N3").Value = -3933.1428  # was -3632.5

# Row 80
$ws.Range("H80").Value = 410.1111  # was 438.25
$ws.Range("I80").Value = 246.66667  # was 277.5
$ws.Range("K80").Value = 246.66667  # was 277.5
$ws.Range("M80").Value = 751.3333299999999  # was 720.5

# Row 83
$ws.Range("H83").Value = 410.1111  # was 438.25
$ws.Range("I83").Value = 246.66667  # was 277.5
$ws.Range("K83").Value = 1233.33335  # was 1387.5
$ws.Range("M83").Value = 3758.66665  # was 3604.5

# Row 134
$ws.Range("H134").Value = 4042.5  # was 4105.3774
$ws.Range("I134").Value = 2099.5715  # was 2140.4412
$ws.Range("K134").Value = 6298.7145  # was 6421.323600000001
$ws.Range("M134").Value = -3763.7145  # was -3886.323600000001

# Row 140
$ws.Range("H140").Value = 249747.5  # was 214778
$ws.Range("J140").Value = 249747.5  # was 214778
$ws.Range("L140").Value = 249747.5  # was 214778
$ws.Range("N140").Value = -260107.5  # was -225138


# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Range("H134").Value = 2000  # was 2500
$ws.Range("J134").Value = 0  # was 3000
$ws.Range("L134").Value = 0  # was 9000
$ws.Range("N134").ClearContents()  # was -14070

# Row 140
$ws.Range("H140").Value = 20000  # was 0
$ws.Range("J140").Value = 20000  # was 0
$ws.Range("L140").Value = 20000  # was 0
$ws.Range("N140").Value = -30360  # was None

# Row 141
$ws.Range("H141").Value = 280463.34  # was 373695
$ws.Range("J141").Value = 280463.34  # was 373695
$ws.Range("L141").Value = 280463.34  # was 373695
$ws.Range("N141").Value = -290823.34  # was -384055


# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 268  # was 218.6
$ws.Range("I2").Value = 0  # was 21
$ws.Range("K2").Value = 0  # was 126
$ws.Range("M2").ClearContents()  # was -13

# Row 37
$ws.Range("H37").Value = 200111970  # was 125122460
$ws.Range("J37").Value = 200111970  # was 125122460
$ws.Range("L37").Value = 600335910  # was 375367380
$ws.Range("N37").Value = -600336134  # was -375367604

# Row 38
$ws.Range("H38").Value = 147.92592  # was 135
$ws.Range("I38").Value = 84.625  # was 80.70587999999999
$ws.Range("J38").Value = 240  # was 206
$ws.Range("K38").Value = 253.875  # was 242.11764
$ws.Range("L38").Value = 720  # was 618
$ws.Range("M38").Value = 93.125  # was 104.88236
$ws.Range("N38").Value = -1414  # was -1312

# Row 107
$ws.Range("H107").Value = 994156.1  # was 1060401.2
$ws.Range("J107").Value = 1589273.8  # was 1765806.4
$ws.Range("L107").Value = 4767821.4  # was 5297419.199999999
$ws.Range("N107").Value = -4771661.4  # was -5301259.199999999

# Row 129
$ws.Range("H129").Value = 2711.4  # was 2226.8667
$ws.Range("I129").Value = 1357  # was 1197.5
$ws.Range("J129").Value = 4065.8  # was 2913.111
$ws.Range("K129").Value = 4071  # was 3592.5
$ws.Range("L129").Value = 12197.4  # was 8739.332999999999
$ws.Range("M129").Value = 929  # was 1407.5
$ws.Range("N129").Value = -22197.4  # was -18739.333

# Row 131
$ws.Range("H131").Value = 1546082  # was 1572711.9
$ws.Range("J131").Value = 1982441.4  # was 2026461.1
$ws.Range("L131").Value = 5947324.199999999  # was 6079383.300000001
$ws.Range("N131").Value = -5957404.199999999  # was -6089463.300000001


# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1663.5  # was 1840.6428
$ws.Range("I22").Value = 1313.7  # was 1415.4445
$ws.Range("J22").Value = 2246.5  # was 2606
$ws.Range("K22").Value = 1313.7  # was 1415.4445
$ws.Range("L22").Value = 2246.5  # was 2606
$ws.Range("M22").Value = -1018.7  # was -1120.4445
$ws.Range("N22").Value = -2836.5  # was -3196

# Row 27
$ws.Range("H27").Value = 1663.5  # was 1840.6428
$ws.Range("I27").Value = 1313.7  # was 1415.4445
$ws.Range("J27").Value = 2246.5  # was 2606
$ws.Range("K27").Value = 1313.7  # was 1415.4445
$ws.Range("L27").Value = 2246.5  # was 2606
$ws.Range("M27").Value = -1206.7  # was -1308.4445
$ws.Range("N27").Value = -2460.5  # was -2820

# Row 122
$ws.Range("H122").Value = 50005028  # was 52636800
$ws.Range("I122").Value = 111113420  # was 125002440
$ws.Range("K122").Value = 333340260  # was 375007320
$ws.Range("M122").Value = -333337810  # was -375004870

# Row 139
$ws.Range("H139").Value = 87141.8  # was 87342.60000000001
$ws.Range("J139").Value = 96427.5  # was 96678.5
$ws.Range("L139").Value = 96427.5  # was 96678.5
$ws.Range("N139").Value = -106707.5  # was -106958.5

# Row 140
$ws.Range("H140").Value = 61447.816  # was 64692.6
$ws.Range("J140").Value = 61447.816  # was 64692.6
$ws.Range("L140").Value = 61447.816  # was 64692.6
$ws.Range("N140").Value = -71807.81599999999  # was -75052.60000000001

# Row 141
$ws.Range("H141").Value = 64950  # was 45475
$ws.Range("J141").Value = 64950  # was 45475
$ws.Range("L141").Value = 64950  # was 45475
$ws.Range("N141").Value = -75310  # was -55835


# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 10421843  # was 11369173
$ws.Range("I132").Value = 4848.278  # was 5246.5757
$ws.Range("J132").Value = 41672828  # was 45460950
$ws.Range("K132").Value = 14544.834  # was 15739.7271
$ws.Range("L132").Value = 125018484  # was 136382850
$ws.Range("M132").Value = -12014.834  # was -13209.7271
$ws.Range("N132").Value = -125023544  # was -136387910

# Row 136
$ws.Range("H136").Value = 8673.51  # was 8772.01
$ws.Range("I136").Value = 757.3333  # was 812.5454999999999
$ws.Range("J136").Value = 9752.987999999999  # was 9755.763999999999
$ws.Range("K136").Value = 2271.9999  # was 2437.6365
$ws.Range("L136").Value = 29258.964  # was 29267.292
$ws.Range("M136").Value = 278.0001000000002  # was 112.3635000000004
$ws.Range("N136").Value = -34358.964  # was -34367.292

